$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$ws.Range('D2').Value = '26.443.72'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').Value = '1.692.47'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('E4').Value = '  +0.60%  '
Set-TextValue 'D5' '219.08'
$ws.Range('E5').Value = '  +1.47%  '
Set-TextValue 'D6' '0.5546'
$ws.Range('E6').Value = '  +8.75%  '
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('E8').Value = '  +1.96%  '
Set-TextValue 'D9' '0.06491'
$ws.Range('E9').Value = '  +1.54%  '
Set-TextValue 'D10' '22.15'
$ws.Range('E10').Value = '  +1.26%  '
Set-TextValue 'D11' '0.07620'
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D12' '4.569'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.687.89'
$ws.Range('E13').Value = '  +0.94%  '
Set-TextValue 'D14' '0.5852'
$ws.Range('E14').Value = '  +0.74%  '
Set-TextValue 'D15' '0.000008491'
$ws.Range('E15').Value = '  -0.27%  '
Set-TextValue 'D16' '65.44'
$ws.Range('E16').Value = '  +2.12%  '
$ws.Range('D17').Value = '26.512.46'
$ws.Range('E17').Value = '  +1.63%  '
Set-TextValue 'D18' '4.975'
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('E19').Value = '  +0.53%  '
$ws.Range('E20').Value = '  +1.94%  '
Set-TextValue 'D21' '191.12'
$ws.Range('E21').Value = '  +0.57%  '
Set-TextValue 'D22' '6.275'
$ws.Range('E22').Value = '  +1.42%  '
Set-TextValue 'D23' '1.011'
$ws.Range('E23').Value = '  +0.56%  '
Set-TextValue 'D24' '150.30'
$ws.Range('E24').Value = '  +3.48%  '
Set-TextValue 'D25' '0.1315'
$ws.Range('E25').Value = '  +8.80%  '
Set-TextValue 'D26' '7.944'
$ws.Range('E26').Value = '  +4.45%  '
Set-TextValue 'D27' '15.83'
$ws.Range('E27').Value = '  +1.23%  '
Set-TextValue 'D28' '0.06353'
$ws.Range('E28').Value = '  -4.43%  '
Set-TextValue 'D29' '1.416'
$ws.Range('E29').Value = '  +6.59%  '
Set-TextValue 'D30' '1.330'
$ws.Range('E30').Value = '  +1.29%  '
Set-TextValue 'D31' '3.597'
$ws.Range('E31').Value = '  +1.36%  '
Set-TextValue 'D32' '3.596'
$ws.Range('E32').Value = '  +2.32%  '
Set-TextValue 'D33' '1.680'
$ws.Range('E33').Value = '  +1.01%  '
Set-TextValue 'D34' '1.049'
$ws.Range('E34').Value = '  +3.28%  '
Set-TextValue 'D35' '0.6268'
$ws.Range('E35').Value = '  +2.21%  '
Set-TextValue 'D36' '2.405'
$ws.Range('E36').Value = '  +1.37%  '
Set-TextValue 'D37' '2.722'
$ws.Range('E37').Value = '  +1.22%  '
Set-TextValue 'D38' '6.257'
$ws.Range('E38').Value = '  -2.01%  '
$ws.Range('D39').Value = '1.127.68'
$ws.Range('E39').Value = '  +3.12%  '
Set-TextValue 'D40' '0.01649'
$ws.Range('E40').Value = '  +3.58%  '
Set-TextValue 'D41' '0.8865'
$ws.Range('E41').Value = '  +1.93%  '
$ws.Range('E42').Value = '  +0.68%  '
Set-TextValue 'D43' '100.67'
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('D44').Value = '1.843.62'
$ws.Range('E44').Value = '  +1.68%  '
Set-TextValue 'D45' '0.00000000109'
$ws.Range('E45').Value = '  -4.55%  '
Set-TextValue 'D46' '57.70'
$ws.Range('E46').Value = '  +2.48%  '
Set-TextValue 'D47' '8.256'
$ws.Range('E47').Value = '  +2.29%  '
$ws.Range('E48').Value = '  +0.17%  '
Set-TextValue 'D49' '0.05283'
$ws.Range('E49').Value = '  +1.09%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D50' '0.4301'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D51' '6.094'
$ws.Range('E51').Value = '  +1.58%  '
